$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the "Periodo Mora" column (E) for every data row (16 through 62,
# the last data row) - matches the new centered style applied across the table.
$ws.Range("E16:E62").HorizontalAlignment = -4108

# Replace the last worker record (row 62): drop JORGE ANDRES MERCADO ZABALETA /
# period 2508 and add a new period (2509) for WILSON ENRIQUE MONTES SUAREZ.
$ws.Range("B62").Value = "CC"
$ws.Range("C62").Value = "1193533559"
$ws.Range("D62").Value = "WILSON ENRIQUE MONTES SUAREZ"
$ws.Range("E62").Value = "2509"
$ws.Range("F62").Value = 100000
$ws.Range("G62").Value = 2500000

# Refresh the summary figures (Cant. Trabajadores, Cant. Periodos, Valor Mora).
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 46
$ws.Range("E11").Value = 3431947
